# Split the "Pearson Correlation (R)" run into four runs:
#   "Pearson Correlation " / "(" / "R" / ")"
# on each of the three labelled textboxes (TextBox 14/15/16) that live
# inside the "Group 17" group shape on slide 4.
#
# Forcing a (no-op) per-character font assignment on the "(" and "R"
# characters causes PowerPoint to break the single run into separate
# runs at those character boundaries, matching the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$grp = $s.Shapes.Item(1)

foreach ($idx in 14, 15, 16) {
    $shp = $grp.GroupItems.Item($idx)
    $tr = $shp.TextFrame.TextRange

    if ($tr.Text -eq "Pearson Correlation (R)") {
        # Characters are 1-based:
        #  1-20 -> "Pearson Correlation "
        #  21   -> "("
        #  22   -> "R"
        #  23   -> ")"
        $tr.Characters(21, 1).Font.Name = "Helvetica"
        $tr.Characters(22, 1).Font.Name = "Helvetica"
        $tr.Characters(23, 1).Font.Name = "Helvetica"
    }
}
